# Append a 4th bullet point to the "Updates by Megh" note in cell D4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d4 = $ws.Range("D4")
$d4.Value2 = $d4.Value2 + "`n4. Deleted folders for the species which are not considered"

# The extra line makes the wrapped text taller, so the row grows accordingly.
$ws.Rows.Item(4).RowHeight = 86.4

# Move the active selection from F4 to D4 (the cell that was edited).
$ws.Range("D4").Select()
